{"js": "// Replace the 25 three-digit \u00f7 one-digit division prompts in the worksheet\n// table with their new values, preserving all existing run/paragraph\n// formatting (font, size, justification, etc).\nconst replacements = [\n  [\"312\u00f72=\", \"814\u00f74=\"],\n  [\"388\u00f72=\", \"173\u00f79=\"],\n  [\"444\u00f75=\", \"412\u00f78=\"],\n  [\"819\u00f76=\", \"624\u00f79=\"],\n  [\"837\u00f73=\", \"705\u00f79=\"],\n  [\"778\u00f78=\", \"570\u00f75=\"],\n  [\"264\u00f75=\", \"520\u00f72=\"],\n  [\"328\u00f78=\", \"738\u00f72=\"],\n  [\"741\u00f79=\", \"740\u00f79=\"],\n  [\"638\u00f78=\", \"879\u00f73=\"],\n  [\"588\u00f76=\", \"391\u00f77=\"],\n  [\"257\u00f79=\", \"992\u00f73=\"],\n  [\"261\u00f72=\", \"290\u00f74=\"],\n  [\"309\u00f74=\", \"842\u00f78=\"],\n  [\"324\u00f75=\", \"965\u00f72=\"],\n  [\"120\u00f79=\", \"586\u00f75=\"],\n  [\"344\u00f79=\", \"400\u00f79=\"],\n  [\"194\u00f72=\", \"141\u00f73=\"],\n  [\"254\u00f74=\", \"764\u00f73=\"],\n  [\"575\u00f74=\", \"175\u00f77=\"],\n  [\"132\u00f73=\", \"458\u00f79=\"],\n  [\"606\u00f77=\", \"284\u00f78=\"],\n  [\"619\u00f77=\", \"598\u00f74=\"],\n  [\"780\u00f77=\", \"432\u00f79=\"],\n  [\"190\u00f75=\", \"287\u00f72=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the 25 three-digit \u00f7 one-digit division prompts in the worksheet\n# table with their new values, preserving all existing formatting.\n$replacements = @(\n    @(\"312\u00f72=\", \"814\u00f74=\"),\n    @(\"388\u00f72=\", \"173\u00f79=\"),\n    @(\"444\u00f75=\", \"412\u00f78=\"),\n    @(\"819\u00f76=\", \"624\u00f79=\"),\n    @(\"837\u00f73=\", \"705\u00f79=\"),\n    @(\"778\u00f78=\", \"570\u00f75=\"),\n    @(\"264\u00f75=\", \"520\u00f72=\"),\n    @(\"328\u00f78=\", \"738\u00f72=\"),\n    @(\"741\u00f79=\", \"740\u00f79=\"),\n    @(\"638\u00f78=\", \"879\u00f73=\"),\n    @(\"588\u00f76=\", \"391\u00f77=\"),\n    @(\"257\u00f79=\", \"992\u00f73=\"),\n    @(\"261\u00f72=\", \"290\u00f74=\"),\n    @(\"309\u00f74=\", \"842\u00f78=\"),\n    @(\"324\u00f75=\", \"965\u00f72=\"),\n    @(\"120\u00f79=\", \"586\u00f75=\"),\n    @(\"344\u00f79=\", \"400\u00f79=\"),\n    @(\"194\u00f72=\", \"141\u00f73=\"),\n    @(\"254\u00f74=\", \"764\u00f73=\"),\n    @(\"575\u00f74=\", \"175\u00f77=\"),\n    @(\"132\u00f73=\", \"458\u00f79=\"),\n    @(\"606\u00f77=\", \"284\u00f78=\"),\n    @(\"619\u00f77=\", \"598\u00f74=\"),\n    @(\"780\u00f77=\", \"432\u00f79=\"),\n    @(\"190\u00f75=\", \"287\u00f72=\")\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
